# Client registration almost done, to-do: transactions
# Add the newly-registered clients below the existing rows on Feuil1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New clients: card number, name (bonus points and credit balance start at 0)
$newClients = @(
    @("953596", "Hana Murata"),
    @("119511", "Jihyeon Nam"),
    @("658433", "Kosuke Yokono"),
    @("917859", "Dahyeon Nam"),
    @("153634", "Celica Puth"),
    @("598450", "Vanita Puth"),
    @("891090", "Melissa Dupuch")
)

$startRow = 23
for ($i = 0; $i -lt $newClients.Count; $i++) {
    $row = $startRow + $i
    $client = $newClients[$i]

    # Card numbers are stored as text (like the existing rows), not numbers.
    # Build them via a TEXT() formula, then paste-special as values so the
    # cell ends up holding a plain text literal (no quote-prefix / number
    # format styling gets attached to the cell).
    $cardCell = $ws.Cells.Item($row, 1)
    $cardCell.Formula = '=TEXT("' + $client[0] + '","0")'
    $cardCell.Copy()
    $cardCell.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = $client[1]
    $ws.Cells.Item($row, 3).Value = 0.0
    $ws.Cells.Item($row, 4).Value = 0.0
}
